$wb = $excel.ActiveWorkbook

# Overview sheet: row 3 is the de-de row ("Ready for handoff")
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-36-11 09:36:13"

# zh-cn sheet: row 3 (d3f61a6d...) status + handoff datetime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-11 09:36:10"

# de-de sheet: row 3 (d3f61a6d...) status + handoff datetime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-11 09:36:13"
